$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 20:11"

# --- Update country case figures (rows identified by current row number) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 3507632
$ws.Range("C4").Value = 28149
$ws.Range("D4").Value = 1565551
$ws.Range("E4").Value = 1803381
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 453
$ws.Range("H4").Value = 138700

# India (row 6)
$ws.Range("B6").Value = 936628
$ws.Range("C6").Value = 28983
$ws.Range("D6").Value = 593080
$ws.Range("E6").Value = 319233
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 588
$ws.Range("H6").Value = 24315

# Turquia (row 18)
$ws.Range("B18").Value = 214993
$ws.Range("C18").Value = 992
$ws.Range("D18").Value = 196720
$ws.Range("E18").Value = 12871
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = 5402

# Canada (row 23)
$ws.Range("B23").Value = 108377
$ws.Range("C23").Value = 222
$ws.Range("D23").Value = 72079
$ws.Range("E23").Value = 27502
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 8796

# Moldavia (row 61)
$ws.Range("B61").Value = 19708
$ws.Range("C61").Value = 269
$ws.Range("D61").Value = 13033
$ws.Range("E61").Value = 6020
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 6
$ws.Range("H61").Value = 655

# Marruecos (row 65)
$ws.Range("B65").Value = 16097
$ws.Range("C65").Value = 161
$ws.Range("D65").Value = 13442
$ws.Range("E65").Value = 2398
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 257

# Guayana Francesa (row 91)
$ws.Range("B91").Value = 6229
$ws.Range("C91").Value = 59
$ws.Range("D91").Value = 3472
$ws.Range("E91").Value = 2726
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = 31

# Maldivas (row 108)
$ws.Range("B108").Value = 2801
$ws.Range("C108").Value = 39
$ws.Range("D108").Value = 2302
$ws.Range("E108").Value = 485
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 14

# Rows 128-130: Libia/Yemen swap places (Yemen now ranks above Libia) and
# Suazilandia's figures are refreshed.
# Row 128 was Libia -> becomes Yemen with updated figures
$ws.Range("A128").Value = "Yemen"
$ws.Range("B128").Value = 1516
$ws.Range("C128").Value = 18
$ws.Range("D128").Value = 685
$ws.Range("E128").Value = 402
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 5
$ws.Range("H128").Value = 429

# Row 129 was Yemen -> becomes Libia, keeping Libia's former figures
$ws.Range("A129").Value = "Libia"
$ws.Range("B129").Value = 1512
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 367
$ws.Range("E129").Value = 1105
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 40

# Row 130 remains Suazilandia, but with refreshed figures
$ws.Range("B130").Value = 1434
$ws.Range("C130").Value = 45
$ws.Range("D130").Value = 695
$ws.Range("E130").Value = 719
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 20

# Republica del Chad (row 146)
$ws.Range("B146").Value = 884
$ws.Range("C146").Value = 4
$ws.Range("D146").Value = 798
$ws.Range("E146").Value = 11
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 75
